$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '40.120.62'
$ws.Range('E2').Value = '  +0.16%  '
$ws.Range('D3').Value = '2.225.39'
$ws.Range('E3').Value = '  -0.76%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '291.50'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.84%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '87.74'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.512'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.59%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.473'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.27%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '30.45'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.51%  '
$ws.Range('E11').Value = '  -2.60%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.47'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.27%  '
$ws.Range('D14').Value = '2.571.64'
$ws.Range('E14').Value = '  -0.58%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '13.93'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.35%  '
$ws.Range('D16').Value = '2.241.60'
$ws.Range('E16').Value = '  -0.07%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.730'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.48%  '
$ws.Range('D18').Value = '40.072.72'
$ws.Range('E18').Value = '  +0.26%  '
$ws.Range('E19').Value = '  -0.98%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.36'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +6.96%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.82'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.33%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '65.69'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.66'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.46'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.28%  '
$ws.Range('E26').Value = '  -0.74%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.75'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.60%  '
$ws.Range('E28').Value = '  -1.51%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.25'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.42%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '156.12'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.91%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '31.86'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.98%  '
$ws.Range('E32').Value = '  +0.01%  '
$ws.Range('E33').Value = '  +1.49%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0719'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.64%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.90'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +6.79%  '
$ws.Range('E36').Value = '  -1.45%  '
$ws.Range('E37').Value = '  +0.19%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '15.77'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.98%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0985'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.97%  '
$ws.Range('E40').Value = '  +0.88%  '
$ws.Range('D41').Value = '2.131.44'
$ws.Range('E41').Value = '  +8.16%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.87'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.54%  '
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '18.31'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +11.90%  '
$ws.Range('B44').Value = 'ApeXProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.14'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.90%  '
$ws.Range('E45').Value = '  -1.25%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.92'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.07%  '
$ws.Range('E47').Value = '  +2.87%  '
$ws.Range('D48').Value = '2.438.09'
$ws.Range('E48').Value = '  -0.81%  '
$ws.Range('E49').Value = '  -0.10%  '
$ws.Range('B50').Value = 'BitcoinSV'
$ws.Range('C50').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '69.48'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.42%  '
$ws.Range('B51').Value = 'TrustWalletToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.10'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.96%  '
